# ============================================================================
# Deutsch-Jozsa Benchmark-Results.xlsx update
#
# Commit: "Adding Benchmarks of Grover's and Amplitude Estimation performed
#          on IONQ simulators"
#
# Sheet 1 "ionq_simulator-0.0.1--ideal":
#   - existing rows 2-18 gain centered alignment (style index 1)
#   - two new benchmark blocks are appended as rows 19-34
#   - dimension grows from A1:S18 to A1:S34
#   - six new merged label rows are added
#
# Sheet 2 "ionq_simulator-0.0.1--aria-2":
#   - trailing blank row 30 is removed (dimension A1:S30 -> A1:S29)
#
# Sheet 3 "ionq_simulator-0.0.1--harmony-1":
#   - trailing blank row 16 is removed (dimension A1:S16 -> A1:S15)
# ============================================================================

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ionq_simulator-0.0.1--ideal")
$ws2 = $wb.Worksheets.Item("ionq_simulator-0.0.1--aria-2")
$ws3 = $wb.Worksheets.Item("ionq_simulator-0.0.1--harmony-1")

# ----------------------------------------------------------------------------
# Sheet 1: center-align the pre-existing rows 2:18 (adds style s="1" to every
# cell in that block, matching the new XML which stamps s="1" throughout).
# ----------------------------------------------------------------------------
$ws1.Range("A2:S18").HorizontalAlignment = -4108

# ----------------------------------------------------------------------------
# Sheet 1: new blank separator row 19 (styled, empty)
# ----------------------------------------------------------------------------
$ws1.Range("A19:S19").HorizontalAlignment = -4108

# ----------------------------------------------------------------------------
# Sheet 1: new benchmark block #1 -> rows 20-25
# ----------------------------------------------------------------------------
$ws1.Range("A20:S22").HorizontalAlignment = -4108
$ws1.Range("A20").Value = "Qiskit-IONQ: Algorithm = Deutsch-Jozsa Simulator = ionq_simulator-0.0.1--ideal"
$ws1.Range("A20:S20").Merge()

$ws1.Range("A21").Value = "CLOUD SIMULATOR - Maximum Supported qubits:29"
$ws1.Range("A21:S21").Merge()

$ws1.Range("A22").Value = "Configuration: Min_Qubits = 3 Max_Qubits = 4 Skip_Qubits = 1 num_circuits = 2  QV_ = None Last_Updated = 2024-09-09 12:05:19"
$ws1.Range("A22:S22").Merge()

$ws1.Range("A23:S23").HorizontalAlignment = -4108
$ws1.Range("A23").Value = "Number of Qubits"
$ws1.Range("B23").Value = "avg_creation_times (ms)"
$ws1.Range("C23").Value = "std_creation_times (ms)"
$ws1.Range("D23").Value = "avg_elapsed_times (ms)"
$ws1.Range("E23").Value = "std_elapsed_times (ms)"
$ws1.Range("F23").Value = "avg_quantum_times (ms)"
$ws1.Range("G23").Value = "std_quantum_times (ms)"
$ws1.Range("H23").Value = "avg_circuit_depths"
$ws1.Range("I23").Value = "avg_transpiled_depths"
$ws1.Range("J23").Value = "Average_Rescaled_fidelity"
$ws1.Range("K23").Value = "Average_Hellinger_fidelity"
$ws1.Range("L23").Value = "std_Rescaled_Fidelity"
$ws1.Range("M23").Value = "std_hellinger_fidelity"
$ws1.Range("N23").Value = "avg_1Q_algorithmic_gate_counts"
$ws1.Range("O23").Value = "avg_2Q_algorithmic_gate_counts"
$ws1.Range("P23").Value = "avg_xi (n2q/n1q+n2q)"
$ws1.Range("Q23").Value = "avg_1Q_Transpiled_gate_counts"
$ws1.Range("R23").Value = "avg_2Q_Transpiled_gate_counts"
$ws1.Range("S23").Value = "avg_tr_xi (tr_n2q/tr_n1q+tr_n2q)"

$ws1.Range("A24:S24").HorizontalAlignment = -4108
$ws1.Range("A24").Value = 3
$ws1.Range("B24").Value = 7.353
$ws1.Range("C24").Value = 2.674
$ws1.Range("D24").Value = 7130.367
$ws1.Range("E24").Value = 96.111
$ws1.Range("F24").Value = 98
$ws1.Range("G24").Value = 40.5
$ws1.Range("H24").Value = 7
$ws1.Range("I24").Value = 7
$ws1.Range("J24").Value = 1
$ws1.Range("K24").Value = 1
$ws1.Range("L24").Value = 0
$ws1.Range("M24").Value = 0
$ws1.Range("N24").Value = 9
$ws1.Range("O24").Value = 1
$ws1.Range("P24").Value = 0.08
$ws1.Range("Q24").Value = 11
$ws1.Range("R24").Value = 1
$ws1.Range("S24").Value = 0.07000000000000001

$ws1.Range("A25:S25").HorizontalAlignment = -4108
$ws1.Range("A25").Value = 4
$ws1.Range("B25").Value = 7.726
$ws1.Range("C25").Value = 0.232
$ws1.Range("D25").Value = 4675.46
$ws1.Range("E25").Value = 1122.592
$ws1.Range("F25").Value = 138
$ws1.Range("G25").Value = 61.5
$ws1.Range("H25").Value = 8
$ws1.Range("I25").Value = 8
$ws1.Range("J25").Value = 1
$ws1.Range("K25").Value = 1
$ws1.Range("L25").Value = 0
$ws1.Range("M25").Value = 0
$ws1.Range("N25").Value = 12.5
$ws1.Range("O25").Value = 1.5
$ws1.Range("P25").Value = 0.09
$ws1.Range("Q25").Value = 16.5
$ws1.Range("R25").Value = 1.5
$ws1.Range("S25").Value = 0.07000000000000001

# ----------------------------------------------------------------------------
# Sheet 1: two blank styled separator rows (26, 27)
# ----------------------------------------------------------------------------
$ws1.Range("A26:S27").HorizontalAlignment = -4108

# ----------------------------------------------------------------------------
# Sheet 1: new benchmark block #2 -> rows 28-33
# ----------------------------------------------------------------------------
$ws1.Range("A28:S30").HorizontalAlignment = -4108
$ws1.Range("A28").Value = "Qiskit-IONQ: Algorithm = Deutsch-Jozsa Simulator = ionq_simulator-0.0.1--ideal"
$ws1.Range("A28:S28").Merge()

$ws1.Range("A29").Value = "CLOUD SIMULATOR - Maximum Supported qubits:29"
$ws1.Range("A29:S29").Merge()

$ws1.Range("A30").Value = "Configuration: Min_Qubits = 3 Max_Qubits = 4 Skip_Qubits = 1 num_circuits = 2  QV_ = None Last_Updated = 2024-09-09 12:11:27"
$ws1.Range("A30:S30").Merge()

# Note: rows 31-34 intentionally keep the *default* (general) alignment -
# the source XML leaves these cells without the s="1" style, unlike the
# equivalent rows (23-27) in block #1.
$ws1.Range("A31").Value = "Number of Qubits"
$ws1.Range("B31").Value = "avg_creation_times (ms)"
$ws1.Range("C31").Value = "std_creation_times (ms)"
$ws1.Range("D31").Value = "avg_elapsed_times (ms)"
$ws1.Range("E31").Value = "std_elapsed_times (ms)"
$ws1.Range("F31").Value = "avg_quantum_times (ms)"
$ws1.Range("G31").Value = "std_quantum_times (ms)"
$ws1.Range("H31").Value = "avg_circuit_depths"
$ws1.Range("I31").Value = "avg_transpiled_depths"
$ws1.Range("J31").Value = "Average_Rescaled_fidelity"
$ws1.Range("K31").Value = "Average_Hellinger_fidelity"
$ws1.Range("L31").Value = "std_Rescaled_Fidelity"
$ws1.Range("M31").Value = "std_hellinger_fidelity"
$ws1.Range("N31").Value = "avg_1Q_algorithmic_gate_counts"
$ws1.Range("O31").Value = "avg_2Q_algorithmic_gate_counts"
$ws1.Range("P31").Value = "avg_xi (n2q/n1q+n2q)"
$ws1.Range("Q31").Value = "avg_1Q_Transpiled_gate_counts"
$ws1.Range("R31").Value = "avg_2Q_Transpiled_gate_counts"
$ws1.Range("S31").Value = "avg_tr_xi (tr_n2q/tr_n1q+tr_n2q)"

$ws1.Range("A32").Value = 3
$ws1.Range("B32").Value = 7.165
$ws1.Range("C32").Value = 0.857
$ws1.Range("D32").Value = 7633.599
$ws1.Range("E32").Value = 293.911
$ws1.Range("F32").Value = 92.5
$ws1.Range("G32").Value = 37.75
$ws1.Range("H32").Value = 7.5
$ws1.Range("I32").Value = 7.5
$ws1.Range("J32").Value = 1
$ws1.Range("K32").Value = 1
$ws1.Range("L32").Value = 0
$ws1.Range("M32").Value = 0
$ws1.Range("N32").Value = 9.5
$ws1.Range("O32").Value = 1
$ws1.Range("P32").Value = 0.08
$ws1.Range("Q32").Value = 11.5
$ws1.Range("R32").Value = 1
$ws1.Range("S32").Value = 0.07000000000000001

$ws1.Range("A33").Value = 4
$ws1.Range("B33").Value = 9.329000000000001
$ws1.Range("C33").Value = 0.771
$ws1.Range("D33").Value = 7473.768
$ws1.Range("E33").Value = 97.643
$ws1.Range("F33").Value = 191.5
$ws1.Range("G33").Value = 87.75
$ws1.Range("H33").Value = 7.5
$ws1.Range("I33").Value = 7.5
$ws1.Range("J33").Value = 1
$ws1.Range("K33").Value = 1
$ws1.Range("L33").Value = 0
$ws1.Range("M33").Value = 0
$ws1.Range("N33").Value = 12
$ws1.Range("O33").Value = 1.5
$ws1.Range("P33").Value = 0.09
$ws1.Range("Q33").Value = 16
$ws1.Range("R33").Value = 1.5
$ws1.Range("S33").Value = 0.07000000000000001

# Row 34 - trailing blank row, unstyled (like rows 31-33)
$ws1.Range("A34").Value = ""

# ----------------------------------------------------------------------------
# Sheet 2: drop the trailing blank row so the used range shrinks back to
# A1:S29 (was A1:S30 with an empty row 30).
# ----------------------------------------------------------------------------
$ws2.Rows(30).Delete()

# ----------------------------------------------------------------------------
# Sheet 3: drop the trailing blank row so the used range shrinks back to
# A1:S15 (was A1:S16 with an empty row 16).
# ----------------------------------------------------------------------------
$ws3.Rows(16).Delete()
